$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$maxRow = $used.Rows.Count

for ($r = 1; $r -le $maxRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text
    if ($val -eq $oldText) {
        $cell.Value = $newText
    }
}
